$d = $word.ActiveDocument

# --- 1. Split the "automatic testing is bad and I should feel bad." sentence
#        so the two occurrences of "bad" become bold runs. ---
$sentence = $d.Content
$sentence.Find.Execute("automatic testing is bad and I should feel bad.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sentenceEnd = $sentence.End

$firstBad = $d.Range($sentence.Start, $sentenceEnd)
$firstBad.Find.Execute("bad", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$firstBad.Bold = 1

$secondBad = $d.Range($firstBad.End, $sentenceEnd)
$secondBad.Find.Execute("bad", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$secondBad.Bold = 1

# --- 2. Update the manual-minus-auto confidence interval in the results table. ---
$d.Content.Find.Execute("manual minus auto: -146.85 [-214.52 to -79.18]", $true, $false, $false, $false, $false, $true, 1, $false, "manual minus auto: -146.85 [-212.88 to -80.81]", 2)
